$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.229.53"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.859.42"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "242.28"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "0.7031"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.3118"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "0.07782"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "24.16"
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("D11").Value = "0.07980"
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("D12").Value = "1.870.16"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D13").Value = "93.82"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "5.162"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "0.6959"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "6.364"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "29.268.15"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "0.000008270"
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("D19").Value = "251.55"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("D20").Value = "2.110.54"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "7.504"
$ws.Range("E23").Value = "  -4.54%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "0.1554"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "8.974"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Value = "159.44"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "4.294"
$ws.Range("E30").Value = "  -2.70%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "1.211"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "0.05253"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("D35").Value = "0.7429"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "1.154"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "2.709"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "0.01866"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "1.241.22"
$ws.Range("E39").Value = "  -3.59%  "
$ws.Range("D40").Value = "2.734"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "6.140"
$ws.Range("E41").Value = "  -6.89%  "
$ws.Range("D42").Value = "110.72"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "0.8946"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "70.63"
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "2.006.86"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").Value = "0.5186"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").Value = "1.781"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "9.441"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "0.4289"
$ws.Range("E51").Value = "  -2.31%  "